$d = $word.ActiveDocument

# --- 1. First paragraph: append two trailing spaces to the existing text,
#        then add three red-colored runs reading
#        "(This is a change - Version for main branch)" ---
$p1 = $d.Paragraphs(1)
$r1 = $p1.Range
# Collapse the range to exclude the paragraph mark (keep just the text run)
$r1.SetRange($r1.Start, $r1.End - 1)
$r1.Text = "This is a Microsoft word document.  "

$enDash = [char]0x2013

$ins1 = $d.Range($r1.End, $r1.End)
$ins1.InsertAfter("(This is a change " + $enDash + " Ve")
$ins1.Font.Color = 255

$ins2 = $d.Range($ins1.End, $ins1.End)
$ins2.InsertAfter("rsion for main branch")
$ins2.Font.Color = 255

$ins3 = $d.Range($ins2.End, $ins2.End)
$ins3.InsertAfter(")")
$ins3.Font.Color = 255

# --- 2. Append a new, empty paragraph at the very end of the body (before
#        the sectPr) shaded with fill color F9F9F9 ---
$endRange = $d.Range($d.Content.End, $d.Content.End)
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/></w:pPr></w:p>'
$endRange.InsertXML($newParaXml)

Write-Host $d.Paragraphs(1).Range.Text
Write-Host ("Paragraph count: " + $d.Paragraphs.Count)
